# Update cryptos list: price (D) and volume/1h change (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D/E columns to text format so numeric-looking strings (e.g. "0.658")
# are stored as literal text rather than being auto-converted to numbers,
# matching the original inlineStr/text representation of these cells.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.901.11"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.262.86"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "0.658"
$ws.Range("E5").Value = "  +4.57%  "
$ws.Range("D6").Value = "233.32"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "63.69"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.450"
$ws.Range("E9").Value = "  +4.59%  "
$ws.Range("D10").Value = "0.0976"
$ws.Range("E10").Value = "  -5.32%  "
$ws.Range("D11").Value = "57.56"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "26.58"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "2.600.47"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "15.62"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "6.14"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("D17").Value = "0.841"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "2.266.52"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "43.838.48"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "0.0₃0979"
$ws.Range("D21").Value = "73.80"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "6.15"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "249.46"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "3.68"
$ws.Range("E25").Value = "  +31.28%  "
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("D27").Value = "2.23"
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("D28").Value = "9.91"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "173.82"
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").Value = "21.89"
$ws.Range("E30").Value = "  +4.07%  "
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").Value = "0.128"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("D34").Value = "4.99"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("D35").Value = "0.0684"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("D37").Value = "3.72"
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("D38").Value = "6.42"
$ws.Range("E38").Value = "  -5.76%  "
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "8.73"
$ws.Range("E42").Value = "  +3.62%  "
$ws.Range("D43").Value = "4.53"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("D44").Value = "17.23"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "98.61"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").Value = "0.0952"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("E48").Value = "  +5.31%  "
$ws.Range("D49").Value = "1.457.34"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -3.43%  "
$ws.Range("E51").Value = "  -1.63%  "

# Restore default cell style (the temporary text NumberFormat above would
# otherwise leave a non-default style applied to these cells).
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
